$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.192.92'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').Value = '2.637.31'
$ws.Range('E3').Value = '  -1.14%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.66'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -0.30%  '
$ws.Range('D9').Value = '2.635.69'
$ws.Range('E9').Value = '  -1.18%  '
$ws.Range('E10').Value = '  +0.39%  '
$ws.Range('E11').Value = '  +1.23%  '
$ws.Range('E12').Value = '  +1.66%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.91'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.14%  '
$ws.Range('D14').Value = '3.120.47'
$ws.Range('E14').Value = '  -1.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000186'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.40%  '
$ws.Range('D16').Value = '72.040.78'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.76'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.39%  '
$ws.Range('D18').Value = '2.653.24'
$ws.Range('E18').Value = '  +0.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.11'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.99'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '374.57'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.07'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.02'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.33'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.72%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.20'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.35'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.23%  '
$ws.Range('D28').Value = '2.774.98'
$ws.Range('E28').Value = '  +1.96%  '
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('E30').Value = '  +1.53%  '
$ws.Range('E31').Value = '  -0.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '489.07'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.31'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.20%  '
$ws.Range('E34').Value = '  -0.15%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '161.27'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.46%  '
$ws.Range('E37').Value = '  +8.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.25'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.92'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.90%  '
$ws.Range('E40').Value = '  -0.18%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('E42').Value = '  -3.09%  '
$ws.Range('E43').Value = '  +1.64%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.85'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.15%  '
$ws.Range('E45').Value = '  -1.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.08'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '150.49'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.64'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.541'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('E50').Value = '  -2.85%  '
$ws.Range('E51').Value = '  +0.83%  '
